$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    4   = 8.737899999999998
    6   = 6.225399999999999
    7   = 5.686199999999998
    8   = 6.3118
    16  = 5.087499999999999
    20  = 9.676599999999985
    21  = 9.235499999999991
    28  = 5.739200000000002
    29  = 5.107700000000005
    30  = 5.864199999999999
    32  = 6.893899999999999
    40  = 9.000999999999994
    46  = 5.429000000000005
    51  = 5.713699999999999
    52  = 5.777400000000003
    57  = 5.020199999999996
    59  = 4.982100000000001
    62  = 5.553399999999995
    66  = 5.768
    73  = 9.053299999999998
    74  = 8.837899999999996
    77  = 9.034000000000004
    92  = 4.987399999999999
    100 = 6.5289
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
